# Append a new forecast row (row 20) to Sheet1, continuing the existing
# yearly time series in columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing date cell (A19) down into A20 so
# the new date cell keeps the same date number format / font / border style.
$ws.Range("A19").Copy($ws.Range("A20"))

# Populate the new row's values.
$ws.Cells.Item(20, 1).Value2 = 45986                 # date_of_forecast (2025-11-25)
$ws.Cells.Item(20, 2).Value2 = 2025                  # y_0
$ws.Cells.Item(20, 3).Value2 = 2.46481303148316      # y_0_forecast
$ws.Cells.Item(20, 4).Value2 = 2026                  # y_1
$ws.Cells.Item(20, 5).Value2 = 2.928470412166684     # y_1_forecast
